$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a plain number must keep their original
# text storage (matches the source inlineStr cells in the diff) - force
# text format before assigning so COM does not silently coerce them to numbers.
$ws.Range("D2").Value = "37.030.65"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "2.062.93"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.67"
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("E6").Value = "  +1.93%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.78"
$ws.Range("E8").Value = "  +10.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.14"
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("E10").Value = "  +1.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0796"
$ws.Range("E11").Value = "  +6.66%  "
$ws.Range("E12").Value = "  +5.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.06"
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("D14").Value = "2.363.79"
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.818"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("E16").Value = "  +3.85%  "
$ws.Range("D17").Value = "2.066.59"
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("D18").Value = "36.992.98"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").Value = "0.0₃0940"
$ws.Range("E19").Value = "  +12.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.46"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.20"
$ws.Range("E21").Value = "  +6.23%  "
$ws.Range("E22").Value = "  +2.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.45"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("E25").Value = "  -4.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.36"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.09"
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.10"
$ws.Range("E28").Value = "  -5.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.01"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.62"
$ws.Range("E31").Value = "  +2.68%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0631"
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.06"
$ws.Range("E33").Value = "  +7.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.38"
$ws.Range("E34").Value = "  +6.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0887"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.28"
$ws.Range("E37").Value = "  -6.91%  "
$ws.Range("E38").Value = "  -4.99%  "
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.104"
$ws.Range("E40").Value = "  +23.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.87"
$ws.Range("E41").Value = "  +9.55%  "
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.14"
$ws.Range("E43").Value = "  -2.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.92"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.13"
$ws.Range("E46").Value = "  +38.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.69"
$ws.Range("E47").Value = "  -51.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.44"
$ws.Range("E48").Value = "  +8.30%  "
$ws.Range("D49").Value = "1.298.98"
$ws.Range("E49").Value = "  -3.14%  "
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.17"
$ws.Range("E51").Value = "  +8.21%  "
